$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change colors in row 2
$ws.Range("C2").Value = "red"
$ws.Range("F2").Value = "blue"

# Update the active cell selection to D8
$ws.Range("D8").Select()
